# Weekly update: insert a new week's row of data (Espinaca, Vega Monumental
# Concepcion) at row 31, pushing the existing rows 31-67 down to 32-68.
#
# Inserting the row keeps all of the previously-recorded observations intact
# (including their formatting) and simply shifts them down by one row, while
# the newly inserted row 31 is populated with this week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 31..67 down to 32..68, inserting a blank row 31.
$ws.Rows.Item(31).Insert()

# Populate the new row 31 with the latest week's record.
$ws.Cells.Item(31, 1).Value2  = 11
$ws.Cells.Item(31, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(31, 3).Value2  = "Bíobío"
$ws.Cells.Item(31, 4).Value2  = 44665
$ws.Cells.Item(31, 5).Value2  = 8
$ws.Cells.Item(31, 6).Value2  = 100112012
$ws.Cells.Item(31, 7).Value2  = "Espinaca"
$ws.Cells.Item(31, 8).Value2  = "Sin especificar"
$ws.Cells.Item(31, 9).Value2  = "Primera"
$ws.Cells.Item(31, 10).Value2 = 100
$ws.Cells.Item(31, 11).Value2 = 6500
$ws.Cells.Item(31, 12).Value2 = 7000
$ws.Cells.Item(31, 13).Value2 = 6750
$ws.Cells.Item(31, 14).Value2 = "`$/cuna 10 kilos"
$ws.Cells.Item(31, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(31, 16).Value2 = 675
$ws.Cells.Item(31, 17).Value2 = 10
$ws.Cells.Item(31, 18).Value2 = "Hortaliza"

# Match the date column's number formatting used by the rest of the sheet.
$ws.Cells.Item(31, 4).NumberFormat = $ws.Cells.Item(32, 4).NumberFormat
